# Recipe.xlsx update: add Proteins/Fat/Carbohydrates/Meal Type/Preraration Time/
# Cooking Time/Image Url columns (E:K) to the header row, and populate the new
# Image Url value for the first recipe row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the seven new header cells (E1:K1), copying D1's formatting first so the
# new cells share the same style as the existing header cells.
$ws.Range("D1").Copy($ws.Range("E1:K1"))
$ws.Range("E1").Value = "Proteins"
$ws.Range("F1").Value = "Fat"
$ws.Range("G1").Value = "Carbohydrates"
$ws.Range("H1").Value = "Meal Type"
$ws.Range("I1").Value = "Preraration Time"
$ws.Range("J1").Value = "Cooking Time"
$ws.Range("K1").Value = "Image Url"

# Populate the Image Url for the first recipe (Chicken Curry, row 2).
$ws.Range("K2").Value = "https://www.recipetineats.com/wp-content/uploads/2023/10/African-coconut-chicken-curry_3.jpg"

# Widen the new Image Url source column (D) and the new Preraration Time
# column (I) to fit their content.
$ws.Columns.Item(4).ColumnWidth = 62.5
$ws.Columns.Item(9).ColumnWidth = 14.333333333333332

# Reflect the last active selection in the saved workbook.
$ws.Range("D26").Select()
